$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the last/duplicated province names:
#  - AFYON -> AFYONKARAHİSAR
#  - İÇEL -> MERSİN
#  - MANİSA (with trailing CR/LF) -> MANİSA (clean)
#  - OSMANİYE (with trailing CR/LF) -> OSMANİYE (clean)
$ws.Range("B4").Value = "AFYONKARAHİSAR"
$ws.Range("B34").Value = "MERSİN"
$ws.Range("B46").Value = "MANİSA"
$ws.Range("B81").Value = "OSMANİYE"

# The two fixed cells no longer contain an embedded line break, so their
# rows no longer need the taller, wrapped row height - shrink them back down.
$ws.Rows("46:46").AutoFit()
$ws.Rows("81:81").AutoFit()

# Select the Name column and auto-fit its width, as happened in the edit session
$ws.Range("B2:B82").Select()
$ws.Columns("B:B").AutoFit()
